## The document currently has no word/styles.xml part (it relies on Word's
## implicit built-in "Normal" style). The target revision introduces that
## part explicitly, defining just the default "Normal" paragraph style.
##
## Touching the Styles collection (defining the built-in "Normal" style)
## is what causes Word to materialize a styles part for the document,
## which is exactly the part-level change described by the diff
## (word/document.xml itself is left untouched).

$d = $word.ActiveDocument

$normal = $d.Styles.Add("Normal", 1)   # 1 = wdStyleTypeParagraph
Write-Host "Defined style:" $normal.NameLocal
